$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 4457.0835
$ws.Range("J70").Value = 3598.5
$ws.Range("L70").Value = 10795.5
$ws.Range("N70").Value = -11335.5

$ws.Range("H73").Value = 4457.0835
$ws.Range("J73").Value = 3598.5
$ws.Range("L73").Value = 10795.5
$ws.Range("N73").Value = -12667.5

$ws.Range("H132").Value = 3226.8076
$ws.Range("I132").Value = 3235.88
$ws.Range("K132").Value = 9707.639999999999
$ws.Range("M132").Value = -7177.639999999999

$ws.Range("H137").Value = 2980.9473
$ws.Range("I137").Value = 2341.2856
$ws.Range("K137").Value = 7023.8568
$ws.Range("M137").Value = -4473.8568

$ws.Range("H138").Value = 4251.5
$ws.Range("J138").Value = 4497.5625
$ws.Range("L138").Value = 13492.6875
$ws.Range("N138").Value = -23772.6875

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 16188.125
$ws.Range("I2").Value = 9998
$ws.Range("K2").Value = 9998
$ws.Range("M2").Value = -9885

$ws.Range("H14").Value = 1425
$ws.Range("J14").Value = 1425
$ws.Range("L14").Value = 1425
$ws.Range("N14").Value = -1775

$ws.Range("H32").Value = 1419.3582
$ws.Range("I32").Value = 1218.4667
$ws.Range("J32").Value = 3141.2856
$ws.Range("K32").Value = 1218.4667
$ws.Range("L32").Value = 3141.2856
$ws.Range("M32").Value = -931.4666999999999
$ws.Range("N32").Value = -3715.2856

$ws.Range("H50").Value = 3479.3572
$ws.Range("I50").Value = 1923.5
$ws.Range("J50").Value = 4646.25
$ws.Range("K50").Value = 1923.5
$ws.Range("L50").Value = 4646.25
$ws.Range("M50").Value = -1209.5
$ws.Range("N50").Value = -6074.25

$ws.Range("H116").Value = 16188.125
$ws.Range("I116").Value = 9998
$ws.Range("K116").Value = 9998
$ws.Range("M116").Value = -7704

$ws.Range("H122").Value = 4587.875
$ws.Range("I122").Value = 3579.2856
$ws.Range("K122").Value = 10737.8568
$ws.Range("M122").Value = -8287.856800000001

$ws.Range("H132").Value = 2399.3704
$ws.Range("I132").Value = 1118.85
$ws.Range("K132").Value = 3356.55
$ws.Range("M132").Value = -826.5499999999997

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 16188.125
$ws.Range("I3").Value = 9998
$ws.Range("K3").Value = 9998
$ws.Range("M3").Value = -9884

$ws.Range("H10").Value = 2254.6
$ws.Range("I10").Value = 1749.6
$ws.Range("J10").Value = 2759.6
$ws.Range("K10").Value = 1749.6
$ws.Range("L10").Value = 2759.6
$ws.Range("M10").Value = -1609.6
$ws.Range("N10").Value = -3039.6

$ws.Range("H134").Value = 3884.3225
$ws.Range("I134").Value = 1532.7273
$ws.Range("K134").Value = 4598.1819
$ws.Range("M134").Value = -2063.1819

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 5484.6523
$ws.Range("I58").Value = 3597.923
$ws.Range("K58").Value = 3597.923
$ws.Range("M58").Value = -3394.923

$ws.Range("H69").Value = 0
$ws.Range("I69").Value = 0
$ws.Range("K69").Value = 0
$ws.Range("M69").ClearContents()

$ws.Range("H72").Value = 0
$ws.Range("I72").Value = 0
$ws.Range("K72").Value = 0
$ws.Range("M72").ClearContents()

$ws.Range("H122").Value = 7833.08
$ws.Range("I122").Value = 4652.3335
$ws.Range("K122").Value = 13957.0005
$ws.Range("M122").Value = -11507.0005

$ws.Range("H134").Value = 2396.6
$ws.Range("I134").Value = 1793.1143
$ws.Range("K134").Value = 5379.3429
$ws.Range("M134").Value = -2844.3429

$ws.Range("H136").Value = 5484.6523
$ws.Range("I136").Value = 3597.923
$ws.Range("K136").Value = 10793.769
$ws.Range("M136").Value = -8243.769

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 369.85715
$ws.Range("I17").Value = 149.75
$ws.Range("K17").Value = 449.25
$ws.Range("M17").Value = -280.25

$ws.Range("H68").Value = 2675.4614
$ws.Range("I68").Value = 992.8889
$ws.Range("J68").Value = 3566.2354
$ws.Range("K68").Value = 2978.6667
$ws.Range("L68").Value = 10698.7062
$ws.Range("M68").Value = -2167.6667
$ws.Range("N68").Value = -12320.7062

$ws.Range("H71").Value = 2675.4614
$ws.Range("I71").Value = 992.8889
$ws.Range("J71").Value = 3566.2354
$ws.Range("K71").Value = 8936.000100000001
$ws.Range("L71").Value = 32096.1186
$ws.Range("M71").Value = -4880.000100000001
$ws.Range("N71").Value = -40208.1186

$ws.Range("H95").Value = 17684.334
$ws.Range("J95").Value = 17013.5
$ws.Range("L95").Value = 51040.5
$ws.Range("N95").Value = -55158.5

$ws.Range("H121").Value = 1264409
$ws.Range("I121").Value = 2245658.8
$ws.Range("J121").Value = 2802.1428
$ws.Range("K121").Value = 6736976.399999999
$ws.Range("L121").Value = 8406.428400000001
$ws.Range("M121").Value = -6735666.399999999
$ws.Range("N121").Value = -11026.4284

$ws.Range("H131").Value = 6275436
$ws.Range("I131").Value = 4445.8
$ws.Range("K131").Value = 13337.4
$ws.Range("M131").Value = -8297.400000000001

$ws.Range("H136").Value = 1324.2222
$ws.Range("I136").Value = 1324.2222
$ws.Range("K136").Value = 3972.6666
$ws.Range("M136").Value = 1127.3334

$ws.Range("H137").Value = 1646.7273
$ws.Range("I137").Value = 1711.6
$ws.Range("J137").Value = 998
$ws.Range("K137").Value = 5134.799999999999
$ws.Range("L137").Value = 2994
$ws.Range("M137").Value = -34.79999999999927
$ws.Range("N137").Value = -13194

$ws.Range("H138").Value = 4426.533
$ws.Range("J138").Value = 8333
$ws.Range("L138").Value = 24999
$ws.Range("N138").Value = -35279

$ws.Range("H139").Value = 3196.4
$ws.Range("I139").Value = 1571.8125
$ws.Range("J139").Value = 6084.5557
$ws.Range("K139").Value = 4715.4375
$ws.Range("L139").Value = 18253.6671
$ws.Range("M139").Value = 424.5625
$ws.Range("N139").Value = -28533.6671

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 6475
$ws.Range("I43").Value = 6950
$ws.Range("J43").Value = 6000
$ws.Range("K43").Value = 6950
$ws.Range("L43").Value = 6000
$ws.Range("M43").Value = -6799
$ws.Range("N43").Value = -6302

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6009.6665
$ws.Range("I7").Value = 6019.6665
$ws.Range("J7").Value = 5999.6665
$ws.Range("K7").Value = 6019.6665
$ws.Range("L7").Value = 5999.6665
$ws.Range("M7").Value = -5907.6665
$ws.Range("N7").Value = -6223.6665

$ws.Range("H22").Value = 8899
$ws.Range("I22").Value = 2631.3333
$ws.Range("K22").Value = 2631.3333
$ws.Range("M22").Value = -2336.3333

$ws.Range("H27").Value = 8899
$ws.Range("I27").Value = 2631.3333
$ws.Range("K27").Value = 2631.3333
$ws.Range("M27").Value = -2524.3333

$ws.Range("H30").Value = 2000
$ws.Range("I30").Value = 2000
$ws.Range("K30").Value = 2000
$ws.Range("M30").Value = -1892

$ws.Range("H68").Value = 5307.769
$ws.Range("I68").Value = 2666.5
$ws.Range("K68").Value = 2666.5
$ws.Range("M68").Value = -1917.5

$ws.Range("H71").Value = 5307.769
$ws.Range("I71").Value = 2666.5
$ws.Range("K71").Value = 13332.5
$ws.Range("M71").Value = -9588.5

$ws.Range("H103").Value = 16834
$ws.Range("J103").Value = 16834
$ws.Range("L103").Value = 16834
$ws.Range("N103").Value = -19178

$ws.Range("H126").Value = 6009.6665
$ws.Range("I126").Value = 6019.6665
$ws.Range("J126").Value = 5999.6665
$ws.Range("K126").Value = 18058.9995
$ws.Range("L126").Value = 17998.9995
$ws.Range("M126").Value = -15588.9995
$ws.Range("N126").Value = -22938.9995

$ws.Range("H141").Value = 79999
$ws.Range("J141").Value = 79999
$ws.Range("L141").Value = 79999
$ws.Range("M141").Value = -90359

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 8000
$ws.Range("J62").Value = 8000
$ws.Range("L62").Value = 8000
$ws.Range("N62").Value = -9248

$ws.Range("H65").Value = 8000
$ws.Range("J65").Value = 8000
$ws.Range("L65").Value = 40000
$ws.Range("N65").Value = -46240

$ws.Range("H81").Value = 7632.8335
$ws.Range("I81").Value = 2949.5
$ws.Range("K81").Value = 5899
$ws.Range("M81").Value = -4838

$ws.Range("H84").Value = 7632.8335
$ws.Range("I84").Value = 2949.5
$ws.Range("K84").Value = 29495
$ws.Range("M84").Value = -24191
